$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1385.5714
$ws.Range("I111").Value = 1116.6666
$ws.Range("J111").Value = 2999
$ws.Range("K111").Value = 3349.9998
$ws.Range("L111").Value = 8997
$ws.Range("M111").Value = -282.9998000000001
$ws.Range("N111").Value = -15131

$ws.Range("H132").Value = 4567230
$ws.Range("I132").Value = 5129099
$ws.Range("J132").Value = 2047.25
$ws.Range("K132").Value = 15387297
$ws.Range("L132").Value = 6141.75
$ws.Range("M132").Value = -15384767
$ws.Range("N132").Value = -11201.75

$ws.Range("H137").Value = 1757.5
$ws.Range("I137").Value = 1664.56
$ws.Range("K137").Value = 4993.68
$ws.Range("M137").Value = -2443.68

$ws.Range("H138").Value = 3597.5334
$ws.Range("J138").Value = 3910.2554
$ws.Range("L138").Value = 11730.7662
$ws.Range("N138").Value = -22010.7662

$ws.Range("H141").Value = 2254.611
$ws.Range("I141").Value = 1799
$ws.Range("K141").Value = 5397
$ws.Range("M141").Value = -217

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3324.5833
$ws.Range("I32").Value = 2799.772
$ws.Range("K32").Value = 2799.772
$ws.Range("M32").Value = -2512.772

$ws.Range("H61").Value = 3078.257
$ws.Range("I61").Value = 2401.3572
$ws.Range("K61").Value = 2401.3572
$ws.Range("M61").Value = -2189.3572

$ws.Range("H74").Value = 5710.683
$ws.Range("I74").Value = 1252.8286
$ws.Range("K74").Value = 1252.8286
$ws.Range("M74").Value = -378.8286000000001

$ws.Range("H77").Value = 5710.683
$ws.Range("I77").Value = 1252.8286
$ws.Range("K77").Value = 6264.143
$ws.Range("M77").Value = -1896.143

$ws.Range("H102").Value = 3950.2727
$ws.Range("I102").Value = 3176
$ws.Range("K102").Value = 3176
$ws.Range("M102").Value = -1554

$ws.Range("H110").Value = 6277.1113
$ws.Range("I110").Value = 6893.7896
$ws.Range("K110").Value = 6893.7896
$ws.Range("M110").Value = -4848.7896

$ws.Range("H132").Value = 1542.5405
$ws.Range("I132").Value = 1165.9688
$ws.Range("K132").Value = 3497.9064
$ws.Range("M132").Value = -967.9064000000003

$ws.Range("H136").Value = 3078.257
$ws.Range("I136").Value = 2401.3572
$ws.Range("K136").Value = 7204.071599999999
$ws.Range("M136").Value = -4654.071599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2550
$ws.Range("I105").Value = 2600
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 2600
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -853
$ws.Range("N105").Value = -5994

$ws.Range("H134").Value = 1822.1163
$ws.Range("I134").Value = 1817.881
$ws.Range("K134").Value = 5453.643
$ws.Range("M134").Value = -2918.643

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 358.91666
$ws.Range("I10").Value = 368.81818
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 368.81818
$ws.Range("L10").Value = 250
$ws.Range("M10").Value = -229.81818
$ws.Range("N10").Value = -528

$ws.Range("H14").Value = 1644.5
$ws.Range("I14").Value = 1089
$ws.Range("K14").Value = 1089
$ws.Range("M14").Value = -919

$ws.Range("H31").Value = 33432.484
$ws.Range("I31").Value = 39826.5
$ws.Range("J31").Value = 9683.286
$ws.Range("K31").Value = 39826.5
$ws.Range("L31").Value = 9683.286
$ws.Range("M31").Value = -39531.5
$ws.Range("N31").Value = -10273.286

$ws.Range("H34").Value = 33432.484
$ws.Range("I34").Value = 39826.5
$ws.Range("J34").Value = 9683.286
$ws.Range("K34").Value = 39826.5
$ws.Range("L34").Value = 9683.286
$ws.Range("M34").Value = -39624.5
$ws.Range("N34").Value = -10087.286

$ws.Range("H107").Value = 967.0526
$ws.Range("I107").Value = 510.1111
$ws.Range("J107").Value = 1378.3
$ws.Range("K107").Value = 510.1111
$ws.Range("L107").Value = 1378.3
$ws.Range("M107").Value = 1409.8889
$ws.Range("N107").Value = -5218.3

$ws.Range("H122").Value = 900.8946999999999
$ws.Range("I122").Value = 900.8946999999999
$ws.Range("K122").Value = 2702.6841
$ws.Range("M122").Value = -252.6840999999999

$ws.Range("H134").Value = 9587.75
$ws.Range("I134").Value = 6344.077
$ws.Range("K134").Value = 19032.231
$ws.Range("M134").Value = -16497.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 374.5

$ws.Range("H132").Value = 1180.6923
$ws.Range("I132").Value = 994.3333
$ws.Range("J132").Value = 1600
$ws.Range("K132").Value = 8948.9997
$ws.Range("L132").Value = 14400
$ws.Range("M132").Value = -6418.9997
$ws.Range("N132").Value = -19460

$ws.Range("H140").Value = 3942.4
$ws.Range("I140").Value = 3942.4
$ws.Range("K140").Value = 11827.2
$ws.Range("M140").Value = -6647.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9584542
$ws.Range("J11").Value = 2417.3333
$ws.Range("L11").Value = 2417.3333
$ws.Range("N11").Value = -2695.3333

$ws.Range("H132").Value = 2700.2334
$ws.Range("I132").Value = 2536.625
$ws.Range("J132").Value = 3354.6667
$ws.Range("K132").Value = 7609.875
$ws.Range("L132").Value = 10064.0001
$ws.Range("M132").Value = -5079.875
$ws.Range("N132").Value = -15124.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 3408.3333
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -888

$ws.Range("H15").Value = 3408.3333
$ws.Range("I15").Value = 1000
$ws.Range("K15").Value = 1000
$ws.Range("M15").Value = -830

$ws.Range("H16").Value = 29415374
$ws.Range("I16").Value = 41669190
$ws.Range("K16").Value = 41669190
$ws.Range("M16").Value = -41669020

$ws.Range("H17").Value = 25250.5
$ws.Range("I17").Value = 25250.5
$ws.Range("K17").Value = 25250.5
$ws.Range("M17").Value = -25080.5

$ws.Range("H136").Value = 3367.182
$ws.Range("I136").Value = 3297.138
$ws.Range("K136").Value = 9891.414000000001
$ws.Range("M136").Value = -7341.414000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 5010000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20480

$ws.Range("H121").Value = 75000
$ws.Range("J121").Value = 75000
$ws.Range("L121").Value = 75000
$ws.Range("N121").Value = -78494

$ws.Range("H132").Value = 3661.125
$ws.Range("I132").Value = 3483.976
$ws.Range("K132").Value = 10451.928
$ws.Range("M132").Value = -7921.928

$ws.Range("H136").Value = 1826.0857
$ws.Range("I136").Value = 1730.7667
$ws.Range("K136").Value = 5192.300099999999
$ws.Range("M136").Value = -2642.300099999999

$ws.Range("H137").Value = 98570.664
$ws.Range("J137").Value = 98570.664
$ws.Range("L137").Value = 98570.664
$ws.Range("N137").Value = -108770.664
